$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match-day dates (Excel serial numbers) for rows 2-19, column A
$dates = @{
    2  = 45787
    3  = 45788
    4  = 45786
    5  = 45785
    6  = 45787
    7  = 45785
    8  = 45788
    9  = 45788
    10 = 45787
    11 = 45787
    12 = 45787
    13 = 45788
    14 = 45787
    15 = 45788
    16 = 45786
    17 = 45787
    18 = 45788
    19 = 45787
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dates[$row]
}
